$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Room Size Distribution")
$ws.Range("B6").Value = 125
$ws.Range("B11").Value = 622
$ws.Range("B12").Value = 1085
$ws.Range("B13").Value = 539
$ws.Range("B14").Value = 852
$ws.Range("B15").Value = 227
$ws.Range("B16").Value = 162
$ws.Range("B17").Value = 335
$ws.Range("B18").Value = 129
$ws.Range("B19").Value = 256
$ws.Range("B20").Value = 116
$ws.Range("B21").Value = 128
$ws.Range("B22").Value = 255
$ws.Range("B23").Value = 94
$ws.Range("B24").Value = 146
$ws.Range("B25").Value = 43
$ws.Range("B26").Value = 36
$ws.Range("B27").Value = 93
$ws.Range("B28").Value = 24
$ws.Range("B29").Value = 55
$ws.Range("B31").Value = 11
$ws.Range("B32").Value = 38

$ws = $wb.Worksheets.Item("Rent Distribution")
$ws.Range("B2").Value = 278
$ws.Range("B3").Value = 119
$ws.Range("B5").Value = 50
$ws.Range("B6").Value = 83
$ws.Range("B7").Value = 109
$ws.Range("B8").Value = 172
$ws.Range("B9").Value = 292
$ws.Range("B10").Value = 442
$ws.Range("B11").Value = 624
$ws.Range("B12").Value = 790
$ws.Range("B13").Value = 819
$ws.Range("B14").Value = 841
$ws.Range("B15").Value = 840
$ws.Range("B16").Value = 572
$ws.Range("B17").Value = 483
$ws.Range("B18").Value = 388
$ws.Range("B19").Value = 294
$ws.Range("B20").Value = 202
$ws.Range("B21").Value = 183
$ws.Range("B22").Value = 104
$ws.Range("B23").Value = 56
$ws.Range("B24").Value = 69
$ws.Range("B25").Value = 32
$ws.Range("B26").Value = 66
$ws.Range("B27").Value = 47
$ws.Range("B29").Value = 22
$ws.Range("B30").Value = 12
$ws.Range("B32").Value = 17

$ws = $wb.Worksheets.Item("Squaremeter Price Distribution")
$ws.Range("B2").Value = 270
$ws.Range("B3").Value = 132
$ws.Range("B4").Value = 56
$ws.Range("B5").Value = 75
$ws.Range("B6").Value = 76
$ws.Range("B7").Value = 135
$ws.Range("B8").Value = 153
$ws.Range("B9").Value = 213
$ws.Range("B10").Value = 300
$ws.Range("B11").Value = 335
$ws.Range("B12").Value = 504
$ws.Range("B13").Value = 524
$ws.Range("B14").Value = 572
$ws.Range("B15").Value = 589
$ws.Range("B16").Value = 425
$ws.Range("B17").Value = 558
$ws.Range("B18").Value = 514
$ws.Range("B19").Value = 421
$ws.Range("B21").Value = 252
$ws.Range("B22").Value = 319
$ws.Range("B23").Value = 221
$ws.Range("B24").Value = 154
$ws.Range("B27").Value = 178
$ws.Range("B29").Value = 84
$ws.Range("B40").Value = 10

